$d = $word.ActiveDocument

$find = "V roku Súhvezdie Perzeus 2022: 16. – 25. januára, 7. – 16. novembra, 6. – 15. decembra"
$replace = "V roku 2022 môžete pozorovať súhvezdie Súhvezdie Perzeus: 16. – 25. januára, 7. – 16. novembra, 6. – 15. decembra"

$r = $d.Content
$r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
